# Auto-generated cell value updates (recalculated price/profit figures)
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (22 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(111, 8).Value = 987
$ws.Cells.Item(111, 9).Value = 992
$ws.Cells.Item(111, 11).Value = 2976
$ws.Cells.Item(111, 13).Value = 91
$ws.Cells.Item(132, 8).Value = 5191.811
$ws.Cells.Item(132, 9).Value = 5919.0938
$ws.Cells.Item(132, 10).Value = 537.2
$ws.Cells.Item(132, 11).Value = 17757.2814
$ws.Cells.Item(132, 12).Value = 1611.6
$ws.Cells.Item(132, 13).Value = -15227.2814
$ws.Cells.Item(132, 14).Value = -6671.6
$ws.Cells.Item(136, 8).Value = 105992.43
$ws.Cells.Item(136, 10).Value = 105992.43
$ws.Cells.Item(136, 12).Value = 105992.43
$ws.Cells.Item(136, 14).Value = -116192.43
$ws.Cells.Item(138, 8).Value = 407485.78
$ws.Cells.Item(138, 9).Value = 3246.818
$ws.Cells.Item(138, 10).Value = 494674.56
$ws.Cells.Item(138, 11).Value = 9740.454000000002
$ws.Cells.Item(138, 12).Value = 1484023.68
$ws.Cells.Item(138, 13).Value = -4600.454000000002
$ws.Cells.Item(138, 14).Value = -1494303.68

# --- Sheet: ARM (35 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 251271.84
$ws.Cells.Item(32, 9).Value = 8085.197
$ws.Cells.Item(32, 10).Value = 1142956.2
$ws.Cells.Item(32, 11).Value = 8085.197
$ws.Cells.Item(32, 12).Value = 1142956.2
$ws.Cells.Item(32, 13).Value = -7798.197
$ws.Cells.Item(32, 14).Value = -1143530.2
$ws.Cells.Item(61, 8).Value = 8546.727999999999
$ws.Cells.Item(61, 9).Value = 6901.4
$ws.Cells.Item(61, 10).Value = 25000
$ws.Cells.Item(61, 11).Value = 6901.4
$ws.Cells.Item(61, 12).Value = 25000
$ws.Cells.Item(61, 13).Value = -6689.4
$ws.Cells.Item(61, 14).Value = -25424
$ws.Cells.Item(110, 8).Value = 2946.0833
$ws.Cells.Item(110, 9).Value = 1723.6316
$ws.Cells.Item(110, 10).Value = 7591.4
$ws.Cells.Item(110, 11).Value = 1723.6316
$ws.Cells.Item(110, 12).Value = 7591.4
$ws.Cells.Item(110, 13).Value = 321.3684000000001
$ws.Cells.Item(110, 14).Value = -11681.4
$ws.Cells.Item(132, 8).Value = 1950.8422
$ws.Cells.Item(132, 9).Value = 1382.24
$ws.Cells.Item(132, 10).Value = 3044.3076
$ws.Cells.Item(132, 11).Value = 4146.72
$ws.Cells.Item(132, 12).Value = 9132.9228
$ws.Cells.Item(132, 13).Value = -1616.72
$ws.Cells.Item(132, 14).Value = -14192.9228
$ws.Cells.Item(136, 8).Value = 8546.727999999999
$ws.Cells.Item(136, 9).Value = 6901.4
$ws.Cells.Item(136, 10).Value = 25000
$ws.Cells.Item(136, 11).Value = 20704.2
$ws.Cells.Item(136, 12).Value = 75000
$ws.Cells.Item(136, 13).Value = -18154.2
$ws.Cells.Item(136, 14).Value = -80100

# --- Sheet: BSM (12 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2873.5715
$ws.Cells.Item(86, 9).Value = 2269.3333
$ws.Cells.Item(86, 11).Value = 2269.3333
$ws.Cells.Item(86, 13).Value = -1146.3333
$ws.Cells.Item(89, 8).Value = 2873.5715
$ws.Cells.Item(89, 9).Value = 2269.3333
$ws.Cells.Item(89, 11).Value = 11346.6665
$ws.Cells.Item(89, 13).Value = -5730.666499999999
$ws.Cells.Item(134, 8).Value = 3192.7334
$ws.Cells.Item(134, 9).Value = 2712.9546
$ws.Cells.Item(134, 11).Value = 8138.8638
$ws.Cells.Item(134, 13).Value = -5603.8638

# --- Sheet: CRP (59 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1542.8572
$ws.Cells.Item(22, 9).Value = 825
$ws.Cells.Item(22, 11).Value = 825
$ws.Cells.Item(22, 13).Value = -475
$ws.Cells.Item(31, 8).Value = 4894.12
$ws.Cells.Item(31, 9).Value = 4646.8677
$ws.Cells.Item(31, 10).Value = 5419.5312
$ws.Cells.Item(31, 11).Value = 4646.8677
$ws.Cells.Item(31, 12).Value = 5419.5312
$ws.Cells.Item(31, 13).Value = -4351.8677
$ws.Cells.Item(31, 14).Value = -6009.5312
$ws.Cells.Item(34, 8).Value = 4894.12
$ws.Cells.Item(34, 9).Value = 4646.8677
$ws.Cells.Item(34, 10).Value = 5419.5312
$ws.Cells.Item(34, 11).Value = 4646.8677
$ws.Cells.Item(34, 12).Value = 5419.5312
$ws.Cells.Item(34, 13).Value = -4444.8677
$ws.Cells.Item(34, 14).Value = -5823.5312
$ws.Cells.Item(50, 8).Value = 88092
$ws.Cells.Item(50, 10).Value = 88092
$ws.Cells.Item(50, 12).Value = 88092
$ws.Cells.Item(50, 14).Value = -89342
$ws.Cells.Item(51, 8).Value = 45049.5
$ws.Cells.Item(51, 10).Value = 80099
$ws.Cells.Item(51, 12).Value = 80099
$ws.Cells.Item(51, 14).Value = -81571
$ws.Cells.Item(58, 8).Value = 3810.7
$ws.Cells.Item(58, 9).Value = 3114.375
$ws.Cells.Item(58, 10).Value = 4606.5
$ws.Cells.Item(58, 11).Value = 3114.375
$ws.Cells.Item(58, 12).Value = 4606.5
$ws.Cells.Item(58, 13).Value = -2911.375
$ws.Cells.Item(58, 14).Value = -5012.5
$ws.Cells.Item(61, 8).Value = 45049.5
$ws.Cells.Item(61, 10).Value = 80099
$ws.Cells.Item(61, 12).Value = 80099
$ws.Cells.Item(61, 14).Value = -80795
$ws.Cells.Item(122, 8).Value = 3156.7437
$ws.Cells.Item(122, 9).Value = 3034.3667
$ws.Cells.Item(122, 10).Value = 3564.6667
$ws.Cells.Item(122, 11).Value = 9103.1001
$ws.Cells.Item(122, 12).Value = 10694.0001
$ws.Cells.Item(122, 13).Value = -6653.1001
$ws.Cells.Item(122, 14).Value = -15594.0001
$ws.Cells.Item(132, 8).Value = 20005136
$ws.Cells.Item(132, 9).Value = 23813354
$ws.Cells.Item(132, 11).Value = 71440062
$ws.Cells.Item(132, 13).Value = -71437532
$ws.Cells.Item(136, 8).Value = 3810.7
$ws.Cells.Item(136, 9).Value = 3114.375
$ws.Cells.Item(136, 10).Value = 4606.5
$ws.Cells.Item(136, 11).Value = 9343.125
$ws.Cells.Item(136, 12).Value = 13819.5
$ws.Cells.Item(136, 13).Value = -6793.125
$ws.Cells.Item(136, 14).Value = -18919.5
$ws.Cells.Item(141, 8).Value = 497012.9
$ws.Cells.Item(141, 10).Value = 497012.9
$ws.Cells.Item(141, 12).Value = 497012.9
$ws.Cells.Item(141, 14).Value = -507372.9

# --- Sheet: CUL (16 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 2908.65
$ws.Cells.Item(34, 10).Value = 3132.353
$ws.Cells.Item(34, 12).Value = 9397.059000000001
$ws.Cells.Item(34, 14).Value = -9565.059000000001
$ws.Cells.Item(39, 8).Value = 12287.375
$ws.Cells.Item(39, 10).Value = 12287.375
$ws.Cells.Item(39, 12).Value = 36862.125
$ws.Cells.Item(39, 14).Value = -37450.125
$ws.Cells.Item(55, 8).Value = 4478
$ws.Cells.Item(55, 10).Value = 4924.875
$ws.Cells.Item(55, 12).Value = 14774.625
$ws.Cells.Item(55, 14).Value = -15128.625
$ws.Cells.Item(132, 8).Value = 3716.7917
$ws.Cells.Item(132, 9).Value = 2940.6
$ws.Cells.Item(132, 11).Value = 26465.4
$ws.Cells.Item(132, 13).Value = -23935.4

# --- Sheet: GSM (12 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 5445
$ws.Cells.Item(97, 9).Value = 890
$ws.Cells.Item(97, 11).Value = 890
$ws.Cells.Item(97, 13).Value = -394
$ws.Cells.Item(113, 8).Value = 6364.722
$ws.Cells.Item(113, 9).Value = 6364.722
$ws.Cells.Item(113, 11).Value = 6364.722
$ws.Cells.Item(113, 13).Value = -4194.722
$ws.Cells.Item(132, 8).Value = 2705.8262
$ws.Cells.Item(132, 9).Value = 2333.5
$ws.Cells.Item(132, 11).Value = 7000.5
$ws.Cells.Item(132, 13).Value = -4470.5

# --- Sheet: LTW (23 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 748.6
$ws.Cells.Item(55, 9).Value = 723.25
$ws.Cells.Item(55, 11).Value = 723.25
$ws.Cells.Item(55, 13).Value = -550.25
$ws.Cells.Item(61, 8).Value = 3768.75
$ws.Cells.Item(61, 9).Value = 2836.5625
$ws.Cells.Item(61, 11).Value = 2836.5625
$ws.Cells.Item(61, 13).Value = -2634.5625
$ws.Cells.Item(113, 8).Value = 3768.75
$ws.Cells.Item(113, 9).Value = 2836.5625
$ws.Cells.Item(113, 11).Value = 2836.5625
$ws.Cells.Item(113, 13).Value = -666.5625
$ws.Cells.Item(132, 8).Value = 3823.158
$ws.Cells.Item(132, 9).Value = 3288
$ws.Cells.Item(132, 11).Value = 9864
$ws.Cells.Item(132, 13).Value = -7334
$ws.Cells.Item(136, 8).Value = 14128.286
$ws.Cells.Item(136, 9).Value = 7966
$ws.Cells.Item(136, 10).Value = 18750
$ws.Cells.Item(136, 11).Value = 23898
$ws.Cells.Item(136, 12).Value = 56250
$ws.Cells.Item(136, 13).Value = -21348
$ws.Cells.Item(136, 14).Value = -61350

# --- Sheet: WVR (41 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(51, 8).Value = 33749.75
$ws.Cells.Item(51, 9).Value = 33333
$ws.Cells.Item(51, 11).Value = 33333
$ws.Cells.Item(51, 13).Value = -32823
$ws.Cells.Item(54, 8).Value = 49458.332
$ws.Cells.Item(54, 10).Value = 49458.332
$ws.Cells.Item(54, 12).Value = 49458.332
$ws.Cells.Item(54, 14).Value = -50498.332
$ws.Cells.Item(100, 8).Value = 83334360
$ws.Cells.Item(100, 9).Value = 1196
$ws.Cells.Item(100, 11).Value = 2392
$ws.Cells.Item(100, 13).Value = -1851
$ws.Cells.Item(107, 8).Value = 1022.25
$ws.Cells.Item(107, 10).Value = 1236.5
$ws.Cells.Item(107, 12).Value = 3709.5
$ws.Cells.Item(107, 14).Value = -7549.5
$ws.Cells.Item(113, 8).Value = 718.26086
$ws.Cells.Item(113, 9).Value = 922.3333
$ws.Cells.Item(113, 11).Value = 2766.9999
$ws.Cells.Item(113, 13).Value = -596.9998999999998
$ws.Cells.Item(126, 8).Value = 2635.4
$ws.Cells.Item(126, 9).Value = 2422.0908
$ws.Cells.Item(126, 10).Value = 4199.6665
$ws.Cells.Item(126, 11).Value = 7266.2724
$ws.Cells.Item(126, 12).Value = 12598.9995
$ws.Cells.Item(126, 13).Value = -4796.2724
$ws.Cells.Item(126, 14).Value = -17538.9995
$ws.Cells.Item(132, 8).Value = 4448178.5
$ws.Cells.Item(132, 9).Value = 5380064.5
$ws.Cells.Item(132, 10).Value = 3799.4614
$ws.Cells.Item(132, 11).Value = 16140193.5
$ws.Cells.Item(132, 12).Value = 11398.3842
$ws.Cells.Item(132, 13).Value = -16137663.5
$ws.Cells.Item(132, 14).Value = -16458.3842
$ws.Cells.Item(136, 8).Value = 2536.682
$ws.Cells.Item(136, 9).Value = 2266.6
$ws.Cells.Item(136, 10).Value = 5237.5
$ws.Cells.Item(136, 11).Value = 6799.799999999999
$ws.Cells.Item(136, 12).Value = 15712.5
$ws.Cells.Item(136, 13).Value = -4249.799999999999
$ws.Cells.Item(136, 14).Value = -20812.5

Write-Host "Applied 220 cell updates across 8 sheets"